# Apply the author's content edits to the "Vehicle ..." mini-deck
# (slides 7-9): retitle the section-break slide and fill in the two
# previously-empty bullet lists with the real use-case / test-stats
# content.

$p = $ppt.ActivePresentation

# Slide 7 - section title: "Vehicle Management" -> "Vehicle Payment software"
$s7 = $p.Slides.Item(7)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Vehicle Payment software"

# Slide 8 - "List of Use Cases Tested": fill in the empty content placeholder
$s8 = $p.Slides.Item(8)
$s8.Shapes.Item(2).TextFrame.TextRange.Text = "Login`rEnter vehicle data`rGenerate Bill`rDisplay vehicle data"

# Slide 9 - "Statistics of Test Cases": replace the placeholder bullets with
# the final write-up of testing results
$s9 = $p.Slides.Item(9)
$s9.Shapes.Item(2).TextFrame.TextRange.Text = "Four functionalities were tested`rLogin does not provide any security and is useless. `rInheritance tree is inverted.`rFunctions fails when input is anything other than provided.`rOther functions are working fine."
